$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 1, shifting all existing data down by one row.
$ws.Rows.Item(1).Insert()

# New header cell in column D, row 1 (bold, same style family as the A-column section headers).
$ws.Range("D1").Value = "HMTL (interface pura)"
$ws.Range("D1").Font.Bold = $true

# New value cell in column D, row 2.
$ws.Range("D2").Value = "ok"

# Column D width to fit the new content (matches bestFit behaviour for the longer header text).
$ws.Columns.Item(4).ColumnWidth = 19.66666666666667

# Update the active selection to D10, per the saved view state.
$ws.Range("D10").Select()
